$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after "Sheet1" and name it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "without EndRow"

# Row 1 — numeric filler values A1:E1
$ws2.Range("A1").Value = 1
$ws2.Range("B1").Value = 2
$ws2.Range("C1").Value = 3
$ws2.Range("D1").Value = 4
$ws2.Range("E1").Value = 5

# Row 2 — numeric filler values A2:E2
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 4
$ws2.Range("E2").Value = 5

# Row 3 — the "#! FINISH" marker (set before the "g" strings below so that
# shared-string indices come out in the same order as the target workbook).
$ws2.Range("A3").Value = "#! FINISH"

# Column I/J on rows 1-2 plus the END_ROW marker on row 1
$ws2.Range("I1").Value = "g"
$ws2.Range("J1").Value = "#! END_ROW"
$ws2.Range("I2").Value = "g"

# Explanation text in B3
$ws2.Range("B3").Value = "<-- this ``#! FINISH`` should be in output, because the line above hasn't ``END_ROW`` that couse finish rendering on out of column limit reached (16384 columns)"

# Make the new sheet the active/visible one, with B3 selected, matching the
# workbook's final view state.
[void]$ws2.Activate()
[void]$ws2.Range("B3").Select()
